$d = $word.ActiveDocument

# Remove the pre-existing "_GoBack" bookmark (it will be recreated in the
# new paragraph with the authors' names, mirroring where Word last left
# the edit cursor).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Insert a new paragraph right after the title paragraph, containing the
# report authors' names (one per line via manual line breaks), followed
# by a fresh "_GoBack" bookmark.
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$namesPara = $d.Paragraphs(2).Range
$xml = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t>Andrew Jacob</w:t></w:r>
<w:r><w:br/><w:t>Joshua Johnson</w:t></w:r>
<w:r><w:br/><w:t>Matthew Jiang</w:t></w:r>
<w:r><w:br/><w:t>Joseph Ryan</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$namesPara.InsertXML($xml) | Out-Null
